$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devices")

# Add header label in G1 (same style as the existing header row, row 7)
$ws.Range("G7").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "AC Loading Details Name"

# Add label in G2 (same style as the data rows below, row 8)
$ws.Range("A8").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("G2").Value = "Signal (AC Units)"

# Update the active selection to G2
$ws.Range("G2").Select()
